$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.700.39"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.702.62"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9952"
$ws.Range("E4").Value = "  -0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.69"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9964"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3966"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4069"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9954"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.515"
$ws.Range("E10").Value = "  +5.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.22"
$ws.Range("E11").Value = "  +9.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08789"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.331"
$ws.Range("E13").Value = "  +9.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.30"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001324"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.506"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.701.46"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.11"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07101"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.48"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.748"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9953"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.677.49"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.978"
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.82"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.139"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.36"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.463"
$ws.Range("E31").Value = "  +24.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.886.91"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").Value = "  -9.04%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08673"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.380"
$ws.Range("E35").Value = "  +19.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.16"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.947"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2724"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.80"
$ws.Range("E39").Value = "  -4.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02768"
$ws.Range("E40").Value = "  +7.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08985"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7649"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7223"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.51"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.168"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9945"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.41"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.323"
$ws.Range("E50").Value = "  +13.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07999"
$ws.Range("E51").Value = "  +0.93%  "
